$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D target cells to Text format so that
# numeric-looking price strings are preserved exactly as text
# (matching the original inlineStr cells), then restore default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.206.49"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.855.80"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "241.24"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "0.6996"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.07783"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "0.3073"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "23.71"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "1.859.93"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "5.107"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "92.10"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "6.550"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "0.000008461"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "29.201.31"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "248.27"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "2.107.81"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "7.544"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "0.1508"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").Value = "161.46"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "8.861"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("D29").Value = "1.551"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "4.208"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").Value = "0.7611"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "1.849"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "1.226.26"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D41").Value = "0.8995"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "109.11"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "5.521"
$ws.Range("E44").Value = "  -10.00%  "
$ws.Range("D45").Value = "2.005.54"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("D47").Value = "65.33"
$ws.Range("E47").Value = "  -7.48%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "9.547"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "1.748"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "7.042"
$ws.Range("E51").Value = "  +0.07%  "

# Restore the default (unstyled) cell format for column D, matching the
# original workbook where these cells carried no explicit style.
$ws.Range("D2:D51").Style = "Normal"
